# Generate Report for Handback
# Swap the b8174788 / a5eaf1b1 rows on every sheet (Overview, zh-cn, de-de),
# refresh the "Handed back" status text + handback datetimes, and widen a
# couple of columns, matching the authored diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item(1)

# Row 2 -> b8174788
$ov.Cells.Item(2,1).Value = "b8174788-9a90-4227-8136-0f93962b431e.md"
$ov.Cells.Item(2,2).Value = "e2e\b8174788-9a90-4227-8136-0f93962b431e.md"
$ov.Cells.Item(2,3).Value = ".md"
$ov.Cells.Item(2,5).Value = "Handed back: not in sync with en-US"
$ov.Cells.Item(2,6).Value = "Handed back: not in sync with en-US"
$ov.Cells.Item(2,7).Value = "2016-09-06 06:45:07"

# Row 3 -> a5eaf1b1
$ov.Cells.Item(3,1).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md"
$ov.Cells.Item(3,2).Value = "e2e\a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md"
$ov.Cells.Item(3,3).Value = ".md"
$ov.Cells.Item(3,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,7).Value = "2016-09-06 06:45:07"

# Hyperlinks on column B: rebuild so the display text follows the swapped
# rows while the underlying link targets (rIds / URLs) stay attached to the
# same row position, matching the authored workbook.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "e2e\b8174788-9a90-4227-8136-0f93962b431e.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "e2e\a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")

# Widen columns E/F slightly (29.98 -> 33.46 char units)
$ov.Columns.Item(5).ColumnWidth = 32.6
$ov.Columns.Item(6).ColumnWidth = 32.6

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item(2)

# Row 2 -> b8174788
$zh.Cells.Item(2,1).Value = "b8174788-9a90-4227-8136-0f93962b431e.md"
$zh.Cells.Item(2,3).Value = "Handed back: not in sync with en-US"
$zh.Cells.Item(2,7).Value = "b8174788-9a90-4227-8136-0f93962b431e.b81b39c77959a143b41bf4b3840b5e702d70060a.zh-cn.xlf"
$zh.Cells.Item(2,9).Value = "b8174788-9a90-4227-8136-0f93962b431e.md"
$zh.Cells.Item(2,10).Value = "b8174788-9a90-4227-8136-0f93962b431e.b81b39c77959a143b41bf4b3840b5e702d70060a.zh-cn.xlf"
$zh.Cells.Item(2,11).Value = "2016-09-06 06:50:36"

# Row 3 -> a5eaf1b1
$zh.Cells.Item(3,1).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md"
$zh.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(3,7).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.7ac40d134b778e577f94f8c45fb09a88968b53d8.zh-cn.xlf"
$zh.Cells.Item(3,9).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md"
$zh.Cells.Item(3,10).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.7ac40d134b778e577f94f8c45fb09a88968b53d8.zh-cn.xlf"
$zh.Cells.Item(3,11).Value = "2016-09-06 06:50:36"

# Hyperlinks on columns A and I
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "b8174788-9a90-4227-8136-0f93962b431e.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/17ed73d302d6ef995738cdd2bd34ba4a3fbbd7a8/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "b8174788-9a90-4227-8136-0f93962b431e.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/17ed73d302d6ef995738cdd2bd34ba4a3fbbd7a8/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")

# Widen column C
$zh.Columns.Item(3).ColumnWidth = 32.6

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item(3)

# Row 2 -> b8174788
$de.Cells.Item(2,1).Value = "b8174788-9a90-4227-8136-0f93962b431e.md"
$de.Cells.Item(2,3).Value = "Handed back: not in sync with en-US"
$de.Cells.Item(2,7).Value = "b8174788-9a90-4227-8136-0f93962b431e.b81b39c77959a143b41bf4b3840b5e702d70060a.de-de.xlf"
$de.Cells.Item(2,9).Value = "b8174788-9a90-4227-8136-0f93962b431e.md"
$de.Cells.Item(2,10).Value = "b8174788-9a90-4227-8136-0f93962b431e.b81b39c77959a143b41bf4b3840b5e702d70060a.de-de.xlf"
$de.Cells.Item(2,11).Value = "2016-09-06 06:50:53"

# Row 3 -> a5eaf1b1
$de.Cells.Item(3,1).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md"
$de.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(3,7).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.7ac40d134b778e577f94f8c45fb09a88968b53d8.de-de.xlf"
$de.Cells.Item(3,9).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md"
$de.Cells.Item(3,10).Value = "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.7ac40d134b778e577f94f8c45fb09a88968b53d8.de-de.xlf"
$de.Cells.Item(3,11).Value = "2016-09-06 06:50:53"

# Hyperlinks on columns A and I
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "b8174788-9a90-4227-8136-0f93962b431e.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4fdfedb8560da5c5b18ad52c7ed59622007814d3/e2e/a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md", "", "", "b8174788-9a90-4227-8136-0f93962b431e.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4fdfedb8560da5c5b18ad52c7ed59622007814d3/e2e/b8174788-9a90-4227-8136-0f93962b431e.md", "", "", "a5eaf1b1-894e-473a-94fa-a2367a38b7ac.md")

# Widen column C
$de.Columns.Item(3).ColumnWidth = 32.6

Write-Output "Report regenerated for handback"
